# Generate Report for Handoff
# Updates the localization-status workbook: flips status from
# "In Translation" to "Ready for handoff" and refreshes the related
# timestamps on the Overview / zh-cn / de-de sheets. Also widens the
# status-related columns that grew to fit the new "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$ws_overview = $wb.Worksheets.Item("Overview")
$ws_zhcn     = $wb.Worksheets.Item("zh-cn")
$ws_dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$ws_overview.Range("E2").Value2 = "Ready for handoff"
$ws_overview.Range("F2").Value2 = "Ready for handoff"
$ws_zhcn.Range("C2").Value2 = "Ready for handoff"
$ws_dede.Range("C2").Value2 = "Ready for handoff"

# --- Timestamps ---
$ws_overview.Range("G2").Value2 = "2016-08-31 20:49:14"
$ws_dede.Range("H2").Value2 = "2016-08-31 20:49:14"
$ws_zhcn.Range("H2").Value2 = "2016-08-31 20:49:09"

# --- Column widths: widen status columns to fit "Ready for handoff" ---
$ws_overview.Columns.Item(5).ColumnWidth = 16.3
$ws_overview.Columns.Item(6).ColumnWidth = 16.3
$ws_zhcn.Columns.Item(3).ColumnWidth = 16.3
$ws_dede.Columns.Item(3).ColumnWidth = 16.3
